# Fixed date errors on resume files
#
# The "Study.com" work-experience entry shows the wrong employment
# dates ("5/2021 - 12/2022"). Correct them to "11/2019 - 2/2021"
# while leaving the identical-looking date range under the "Ezoic"
# entry untouched.
#
# The date text is split across two runs in the original document:
#   Run A (bold, normal size):      "5/2021 - 12/202"
#   Run B (bold, slightly larger):  "2"
# so the edit is done as two scoped, offset-based Range replacements
# that line up with those run boundaries instead of a single
# document-wide Find/Replace (which would also touch the Ezoic entry
# and would collapse the two differently-formatted runs into one).

$d = $word.ActiveDocument

# Locate the "Study.com" entry's heading paragraph specifically (not
# the look-alike "Ezoic" entry that shares the exact same date text).
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Study.com*" -and $t -like "*5/2021 - 12/2022*") {
        $target = $p
    }
}

if ($target -ne $null) {
    $pStart = $target.Range.Start
    $pText = $target.Range.Text

    $oldRunA = "5/2021 - 12/202"
    $idxA = $pText.IndexOf($oldRunA)

    if ($idxA -ge 0) {
        # Run A: "5/2021 - 12/202" -> "11/2019 - 2/202"
        $runAStart = $pStart + $idxA
        $runAEnd = $runAStart + $oldRunA.Length
        $runA = $d.Range($runAStart, $runAEnd)
        $runA.Text = "11/2019 - 2/202"

        # Run B immediately follows Run A and holds the trailing digit
        # of the year ("2" -> "1"); its own (larger) font size must be
        # preserved, so it is edited as its own Range rather than as
        # part of a single combined replacement.
        $runBStart = $runAEnd
        $runBEnd = $runBStart + 1
        $runB = $d.Range($runBStart, $runBEnd)
        $runB.Text = "1"
    }
}
